$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (far outside the used range) used to coerce numeric-looking
# strings (e.g. "1.026", "0.07130") into text cells without Excel's
# auto-number-detection turning them into real numbers (which would also
# silently trim insignificant trailing zeros, e.g. "0.07130" -> "0.0713").
# We append a non-numeric sentinel character so Excel stores the literal as
# text, strip the sentinel back off with a formula (result is still text),
# then copy/paste-special just the *value* of that formula cell into the
# real target cell and clear the helpers. This keeps the destination cell
# a plain text cell with no style/number-format side effects.
$helper = $ws.Range("ZZ1")
$helperCalc = $ws.Range("ZZ2")

$ws.Range('D2').Value = '27.408.32'
$ws.Range('E2').Value = '  +3.82%  '
$ws.Range('D3').Value = '1.836.18'
$ws.Range('E3').Value = '  +3.33%  '
$helper.Value = '1.026x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('E4').Value = '  +2.41%  '
$helper.Value = '318.22x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +4.08%  '
$ws.Range('E6').Value = '  +2.22%  '
$ws.Range('E7').Value = '  +2.77%  '
$ws.Range('E8').Value = '  +3.07%  '
$helper.Value = '0.07368x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  +3.21%  '
$helper.Value = '0.8737x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  +4.27%  '
$helper.Value = '21.33x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  +4.36%  '
$ws.Range('D12').Value = '1.891.94'
$ws.Range('E12').Value = '  +6.06%  '
$helper.Value = '5.459x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Value = '  +4.26%  '
$helper.Value = '6.688x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  +3.78%  '
$helper.Value = '0.07130x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  +3.71%  '
$helper.Value = '82.41x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  +4.50%  '
$helper.Value = '1.029x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  +2.32%  '
$helper.Value = '0.000009005x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  +3.75%  '
$helper.Value = '1.023x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  +2.09%  '
$helper.Value = '15.32x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  +2.44%  '
$ws.Range('D21').Value = '27.413.83'
$ws.Range('E21').Value = '  +3.83%  '
$helper.Value = '5.231x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  +3.13%  '
$ws.Range('E23').Value = '  +1.17%  '
$helper.Value = '156.64x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +3.00%  '
$helper.Value = '1.898x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +5.92%  '
$helper.Value = '18.59x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  +3.11%  '
$helper.Value = '5.233x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  +3.33%  '
$helper.Value = '1.918x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  +5.96%  '
$helper.Value = '115.83x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  +1.43%  '
$helper.Value = '0.09039x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  +2.17%  '
$ws.Range('E31').Value = '  +6.96%  '
$helper.Value = '0.7592x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  +4.97%  '
$helper.Value = '4.472x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  +3.34%  '
$helper.Value = '2.862x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  +4.40%  '
$ws.Range('E35').Value = '  +2.43%  '
$ws.Range('E36').Value = '  +4.73%  '
$helper.Value = '0.01963x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  +4.41%  '
$helper.Value = '0.05245x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  +2.21%  '
$helper.Value = '0.5160x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  +5.10%  '
$helper.Value = '2.770x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  +6.02%  '
$ws.Range('E41').Value = '  +3.09%  '
$helper.Value = '6.563x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  +3.60%  '
$helper.Value = '8.486x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  +6.41%  '
$helper.Value = '108.60x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +3.68%  '
$helper.Value = '10.58x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  +3.73%  '
$helper.Value = '1.026x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +2.51%  '
$helper.Value = '1.681x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  +2.36%  '
$helper.Value = '0.4622x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  +4.25%  '
$helper.Value = '1.898x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  +10.59%  '
$helper.Value = '0.06302x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  +1.98%  '
$helper.Value = '39.35x'
$helperCalc.Formula = '=LEFT(ZZ1,LEN(ZZ1)-1)'
$helperCalc.Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +7.18%  '

$helper.Clear()
$helperCalc.Clear()
$excel.CutCopyMode = 0
